$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking data refresh: update Price (D) and Volume(1h) (E) columns.
# D-column text values that look like plain numbers (e.g. "143.52") get forced
# to stay text (matching the source's inlineStr cells) via NumberFormat="@",
# then the format is reset back to Normal so no stray style is left on the cell.
$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = "63.962.41"
$dCell.Style = "Normal"
$ws.Range("E2").Value = "  +1.41%  "
$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = "3.319.50"
$dCell.Style = "Normal"
$ws.Range("E3").Value = "  +6.23%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  +1.16%  "
$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = "143.52"
$dCell.Style = "Normal"
$ws.Range("E6").Value = "  +4.97%  "
$ws.Range("E7").Value = "  -0.01%  "
$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = "3.319.56"
$dCell.Style = "Normal"
$ws.Range("E8").Value = "  +6.49%  "
$ws.Range("E9").Value = "  +1.53%  "
$ws.Range("E10").Value = "  +3.29%  "
$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = "5.57"
$dCell.Style = "Normal"
$ws.Range("E11").Value = "  +5.16%  "
$ws.Range("E12").Value = "  +4.33%  "
$ws.Range("E13").Value = "  +1.50%  "
$ws.Range("E14").Value = "  +2.56%  "
$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = "3.866.41"
$dCell.Style = "Normal"
$ws.Range("E15").Value = "  +6.34%  "
$ws.Range("E16").Value = "  +0.89%  "
$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = "3.318.45"
$dCell.Style = "Normal"
$ws.Range("E17").Value = "  +6.30%  "
$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = "64.044.55"
$dCell.Style = "Normal"
$ws.Range("E18").Value = "  +1.53%  "
$ws.Range("E19").Value = "  +3.96%  "
$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = "482.95"
$dCell.Style = "Normal"
$ws.Range("E20").Value = "  +2.09%  "
$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = "14.30"
$dCell.Style = "Normal"
$ws.Range("E21").Value = "  +1.30%  "
$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = "0.739"
$dCell.Style = "Normal"
$ws.Range("E22").Value = "  +6.20%  "
$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = "8.01"
$dCell.Style = "Normal"
$ws.Range("E23").Value = "  +3.50%  "
$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = "13.66"
$dCell.Style = "Normal"
$ws.Range("E24").Value = "  +5.83%  "
$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = "84.88"
$dCell.Style = "Normal"
$ws.Range("E25").Value = "  -1.62%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  +2.17%  "
$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = "7.43"
$dCell.Style = "Normal"
$ws.Range("E28").Value = "  +6.42%  "
$ws.Range("E29").Value = "  -0.08%  "
$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = "8.21"
$dCell.Style = "Normal"
$ws.Range("E30").Value = "  +3.61%  "
$ws.Range("E31").Value = "  +5.73%  "
$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = "29.57"
$dCell.Style = "Normal"
$ws.Range("E32").Value = "  +11.08%  "
$ws.Range("E33").Value = "  -1.15%  "
$ws.Range("E34").Value = "  +1.78%  "
$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = "1.11"
$dCell.Style = "Normal"
$ws.Range("E35").Value = "  +2.73%  "
$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = "6.01"
$dCell.Style = "Normal"
$ws.Range("E36").Value = "  +3.79%  "
$ws.Range("E37").Value = "  +7.96%  "
$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = "52.86"
$dCell.Style = "Normal"
$ws.Range("E38").Value = "  +1.41%  "
$dCell = $ws.Range("D39")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0405"
$dCell.Style = "Normal"
$ws.Range("E39").Value = "  +4.96%  "
$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = "433.66"
$dCell.Style = "Normal"
$ws.Range("E40").Value = "  +2.51%  "
$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = "3.064.36"
$dCell.Style = "Normal"
$ws.Range("E41").Value = "  +5.69%  "
$ws.Range("E42").Value = "  +3.15%  "
$ws.Range("E43").Value = "  +3.42%  "
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("E45").Value = "  +1.64%  "
$ws.Range("E46").Value = "  +4.88%  "
$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = "26.61"
$dCell.Style = "Normal"
$ws.Range("E47").Value = "  +4.45%  "
$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = "36.06"
$dCell.Style = "Normal"
$ws.Range("E48").Value = "  +15.21%  "
$ws.Range("E50").Value = "  +2.74%  "
$ws.Range("E51").Value = "  +2.11%  "
